$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. "587.90"); Excel would
# otherwise coerce them to real numbers. Force text storage, then restore the
# default "Normal" cell style so no visible formatting change is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '72.082.05'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.36%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.634.80'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.92%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.58%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.520'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.33%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.633.82'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.97%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.171'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.66%  '

$ws.Range("E11").Value = '  +0.94%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.360'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.90'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.38%  '

$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000191'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.81%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.118.26'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.23%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '71.900.68'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.13%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.74'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.67%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.624.03'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.62%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.13'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.39%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.69%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '374.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.91%  '

$ws.Range("E22").Value = '  -1.68%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.66%  '

$ws.Range("E25").Value = '  -0.09%  '

$ws.Range("E26").Value = '  -2.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.94%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.773.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.92%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.987'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.52%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0952'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.93'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.02%  '

$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.33'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.64%  '

$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '487.52'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.32%  '

$ws.Range("E34").Value = '  -0.44%  '

$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.26'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.50%  '

$ws.Range("E37").Value = '  +7.69%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.28'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.92%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.91'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.89%  '

$ws.Range("E40").Value = '  -0.75%  '

$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.72'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.53%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.57'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.62%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.327'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.03'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.34%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '150.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.45%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.65'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.58%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.540'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.24%  '

$ws.Range("E50").Value = '  -3.64%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.606'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.26%  '
